$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift the timesheet one week forward: day labels in column A (rows 11-17)
# go from "30/03 - 05/04" to "06/04 - 12/04".
$ws.Range("A11").Value = "Sun 06/04"
$ws.Range("A12").Value = "Mon 07/04"
$ws.Range("A13").Value = "Tue 08/04"
$ws.Range("A14").Value = "Wed 09/04"
$ws.Range("A15").Value = "Thur 10/04"
$ws.Range("A16").Value = "Fri   11/04"
$ws.Range("A17").Value = "Sat  12/04"

# "Week of:" date (G8) moves back one week to the new week's start date
$ws.Range("G8").Value = 41735

# Refresh the view: scroll position and selection
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$null = $ws.Range("A17").Select()
